# Auto-commit style refresh of the "Metrics" sheet's raw input numbers.
# The "today" sheet pulls these via =Metrics!Bxx formulas (plus a couple of
# derived E/F columns and a TODAY()-1 cell), so those recalc on their own
# once the source values below are updated - no need to touch them directly.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 377126.62
$wsMetrics.Range("B3").Value  = 309332.12000000005
$wsMetrics.Range("B4").Value  = 120213.18
$wsMetrics.Range("B5").Value  = 15053
$wsMetrics.Range("B6").Value  = 4744258.09
$wsMetrics.Range("B7").Value  = 3999150.7899999991
$wsMetrics.Range("B8").Value  = 1390815.3199999998
$wsMetrics.Range("B9").Value  = 184054
$wsMetrics.Range("B10").Value = 33209581.890999824
$wsMetrics.Range("B11").Value = 31274372.310000002
$wsMetrics.Range("B12").Value = 11672524.210000005
$wsMetrics.Range("B13").Value = 1281681

# Restore the recorded cursor position on Metrics (this sheet is not the
# active tab, but Excel still remembers its own last selection).
$wsMetrics.Range("E20").Select()

# Move back to "today" - it was (and stays) the active/visible sheet - and
# leave its cursor where it was last recorded.
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("D8").Select()
